$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text so numeric-looking strings
# (e.g. "2.800", "28.246.48") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.246.48'
$ws.Range('E2').Value = '  +3.67%  '
$ws.Range('D3').Value = '1.917.56'
$ws.Range('E3').Value = '  +3.03%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -1.25%  '
$ws.Range('D5').Value = '316.17'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '0.4872'
$ws.Range('E7').Value = '  +1.67%  '
$ws.Range('D8').Value = '0.3852'
$ws.Range('E8').Value = '  +3.41%  '
$ws.Range('D9').Value = '0.07438'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').Value = '0.9579'
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('D11').Value = '21.18'
$ws.Range('E11').Value = '  +3.82%  '
$ws.Range('D12').Value = '0.07840'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '1.901.62'
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').Value = '5.585'
$ws.Range('E14').Value = '  +2.98%  '
$ws.Range('D15').Value = '6.679'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').Value = '92.58'
$ws.Range('E16').Value = '  +2.73%  '
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').Value = '0.000008938'
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('D20').Value = '28.227.44'
$ws.Range('E20').Value = '  +3.43%  '
$ws.Range('D21').Value = '15.10'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').Value = '5.188'
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('D23').Value = '2.151.91'
$ws.Range('E23').Value = '  +2.14%  '
$ws.Range('E24').Value = '  +2.86%  '
$ws.Range('D25').Value = '1.963'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').Value = '157.28'
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('D27').Value = '18.74'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('D28').Value = '2.126'
$ws.Range('E28').Value = '  +6.29%  '
$ws.Range('D29').Value = '116.99'
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('D30').Value = '5.047'
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('D31').Value = '0.08926'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E33').Value = '  +5.38%  '
$ws.Range('D34').Value = '0.7849'
$ws.Range('E34').Value = '  +5.89%  '
$ws.Range('D35').Value = '4.725'
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('D36').Value = '2.800'
$ws.Range('E36').Value = '  +4.60%  '
$ws.Range('D37').Value = '1.139'
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('D38').Value = '0.02062'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05422'
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5619'
$ws.Range('E40').Value = '  +4.83%  '
$ws.Range('D41').Value = '3.036'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').Value = '7.130'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = '8.662'
$ws.Range('E43').Value = '  +3.76%  '
$ws.Range('D44').Value = '0.1543'
$ws.Range('D45').Value = '0.4979'
$ws.Range('E45').Value = '  +3.64%  '
$ws.Range('D46').Value = '10.83'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').Value = '107.11'
$ws.Range('E47').Value = '  +4.30%  '
$ws.Range('D48').Value = '1.693'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').Value = '1.007'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').Value = '69.66'
$ws.Range('E50').Value = '  +4.49%  '
$ws.Range('D51').Value = '0.06135'
$ws.Range('E51').Value = '  +0.87%  '

# Restore the original (default) style on the Price column now that
# the values have been written as text, so cell styling matches the source.
$ws.Range("D2:D51").Style = "Normal"
